# Update the sample flight date in the emergency checklist form.
# The "Start time" value in row 2 (B2) moves from 2019-01-01 to 2019-01-10,
# keeping the same time-of-day fraction (43466.65121527778 -> 43475.65121527778).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form1")
$ws.Activate()

$ws.Range("B2").Value = 43475.65121527778

# Reflect the user's click on the edited cell as the new selection/view state.
$ws.Range("B2").Select()
